# CA04/grades.xlsx - "adding the grades and fixing a bug in grade.py"
#
# The sheet tracks two quiz attempts (Q1 columns C:E, Q2 columns F:H) for a
# list of students. Column C/F hold the "test case (40)" score. Many of
# those cells were still showing the "**" placeholder (not-yet-graded)
# shared string. This change fills in the actual test-case grades (40/40)
# for the students who have already had their other scores (columns D/E,
# G/H) entered, while students that are still fully ungraded (rows 22-32,
# which still show "**" across every score column) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose Q1 "test case" grade (column C) is being recorded as 40/40.
# (Row 5 already had a numeric grade in C and is skipped.)
$rowsC = @(3, 4, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21)

# Rows whose Q2 "test case" grade (column F) is being recorded as 40/40.
$rowsF = @(6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21)

foreach ($r in $rowsC) {
    # The neighbouring "presentation" cell (column D) already carries the
    # normal graded-cell look (font/fill/border) that a filled-in test-case
    # score should use instead of the placeholder's styling, so copy its
    # format across before writing the real number.
    $ws.Range("D$r").Copy()
    $ws.Range("C$r").PasteSpecial(-4122)
    $ws.Range("C$r").Value = 40
}

foreach ($r in $rowsF) {
    $ws.Range("F$r").Value = 40
}

$excel.CutCopyMode = 0

# Move the active selection to reflect where the author was working.
$ws.Range("F23").Select()
